$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: new Bloomberg alert-log entry (1033 HK, RECOM_SMALL) ---
$ws.Range("A21").Value = "'04/19/2018"
$ws.Range("A21").ClearFormats()
$ws.Range("B21").Value = "'06:16:02"
$ws.Range("B21").ClearFormats()
$ws.Range("C21").Value = "1033 HK Equity"
$ws.Range("D21").Value = "1033 HK"
$ws.Range("E21").Value = " Target Px decreased to 2.913 by J.P. Morgan"
$ws.Range("G21").Value = "RECOM_SMALL"

# --- Row 22: second alert-log entry, same ticker, different time/text ---
$ws.Range("A22").Value = "'04/19/2018"
$ws.Range("A22").ClearFormats()
$ws.Range("B22").Value = "'06:14:02"
$ws.Range("B22").ClearFormats()
$ws.Range("C22").Value = "1033 HK Equity"
$ws.Range("D22").Value = "1033 HK"
$ws.Range("E22").Value = " Target Px decreased to .94 by J.P. Morgan"
$ws.Range("G22").Value = "RECOM_SMALL"

# --- Row 24: new exception row (SINOPEC OILFIE-H, TP+ by J.P. Morgan) ---
$ws.Range("A24").Value = 60
$ws.Range("B24").Value = "B/D"
$ws.Range("B24").ClearFormats()
$ws.Range("C24").Value = "TP+ "
$ws.Range("D2").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("D24").Value = 43209
$ws.Range("E24").Value = "J.P. Morgan"
$ws.Range("F24").Value = "SINOPEC OILFIE-H"
$ws.Range("H24").Value = 157.78761061946901
$ws.Range("I24").Value = -9.7345132743362797
$ws.Range("J2").Copy()
$ws.Range("J24").PasteSpecial(-4122)
$ws.Range("J24").Value = 185.58823529411799
$ws.Range("K24").Value = 131.55802563264999
$ws.Range("L24").Value = 4.1868916
$ws.Range("M24").Value = 3.0678559999999999
$ws.Range("O24").Value = 1
$ws.Range("P24").Value = 1
$ws.Range("R24").Value = 0.5

# --- Row 25: paired exception row (SINOPEC OILFIE-H, TP- by J.P. Morgan) ---
$ws.Range("A25").Value = 61
$ws.Range("B25").Value = "B/D"
$ws.Range("B25").ClearFormats()
$ws.Range("C25").Value = "TP- "
$ws.Range("D2").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D25").Value = 43209
$ws.Range("E25").Value = "J.P. Morgan"
$ws.Range("F25").Value = "SINOPEC OILFIE-H"
$ws.Range("H25").Value = -16.814159292035399
$ws.Range("I25").Value = -9.7345132743362797
$ws.Range("J2").Copy()
$ws.Range("J25").PasteSpecial(-4122)
$ws.Range("J25").Value = -7.8431372549019702
$ws.Range("K25").Value = -25.278220358842699
$ws.Range("L25").Value = 4.1868916
$ws.Range("M25").Value = 3.0678559999999999
$ws.Range("O25").Value = 1
$ws.Range("P25").Value = 1
$ws.Range("R25").Value = -0.5

$excel.CutCopyMode = 0

# --- Column F width (now has data) ---
$ws.Columns("F").ColumnWidth = 17.5703125

# --- Cell comment on F21 (mistake note) ---
$ws.Range("F21").AddComment("YChen:" + [char]10 + "mistake on Bloomberg Alert")

# --- Selection moved to H22 ---
$null = $ws.Range("H22").Select()
